# Update cryptos list values (price + 1h volume change) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.276.06"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.868.15"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4698"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2856"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06565"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07818"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").Value = "1.868.21"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6933"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.079"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "268.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "30.369.98"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007696"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.81%  "
$ws.Range("D21").Value = "2.124.14"
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.255"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.566"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.938"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("E29").Value = "  -2.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09897"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.352"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.054"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("E34").Value = "  +0.81%  "
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7041"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01871"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("E39").Value = "  +5.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.315"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4170"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8377"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.82"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "980.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.121"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.156"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("E51").Value = "  +0.43%  "
